# Generate Report for handback
# Update the "Correspond Handoff Datetime" and "Correspond Handback DateTime"
# timestamps on the per-language sheets (zh-cn, de-de) to reflect the new
# handback report generation times.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D2").Value = "2016-01-17 14:32:08"
$wsZhCn.Range("G2").Value = "2016-01-17 14:32:55"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D2").Value = "2016-01-17 14:32:18"
$wsDeDe.Range("G2").Value = "2016-01-17 14:33:13"
